$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values updated per the refreshed cryptos snapshot.
# Cells whose new text is numeric-looking get an apostrophe prefix so Excel
# keeps storing them as text (matching the original inlineStr cells), then
# their style is reset to Normal so no stray number-format style sticks.

$ws.Range("D2").Value = '25.725.55'
$ws.Range("E2").Value = '  -0.16%  '
$ws.Range("D3").Value = '1.629.95'
$ws.Range("E3").Value = '  -0.34%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = "'214.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.57%  '
$ws.Range("D6").Value = "'0.501"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.56%  '
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("E8").Value = '  -1.07%  '
$ws.Range("E9").Value = '  -1.53%  '
$ws.Range("D10").Value = "'19.49"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.96%  '
$ws.Range("D11").Value = "'0.0791"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.42%  '
$ws.Range("E12").Value = '  +0.08%  '
$ws.Range("D13").Value = '1.855.03'
$ws.Range("E13").Value = '  -0.33%  '
$ws.Range("D14").Value = '1.630.68'
$ws.Range("E14").Value = '  -0.29%  '
$ws.Range("E15").Value = '  +0.05%  '
$ws.Range("D16").Value = '0.0₃0760'
$ws.Range("E16").Value = '  -2.06%  '
$ws.Range("D17").Value = "'63.01"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = '25.737.70'
$ws.Range("E18").Value = '  -0.25%  '
$ws.Range("E19").Value = '  -0.01%  '
$ws.Range("E20").Value = '  -0.58%  '
$ws.Range("D21").Value = "'191.38"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.44%  '
$ws.Range("E22").Value = '  -0.52%  '
$ws.Range("D23").Value = "'6.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.52%  '
$ws.Range("E24").Value = '  -0.05%  '
$ws.Range("E25").Value = '  +3.19%  '
$ws.Range("D26").Value = "'142.96"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.81%  '
$ws.Range("E27").Value = '  +2.51%  '
$ws.Range("D28").Value = "'6.86"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.13%  '
$ws.Range("D29").Value = "'15.45"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.80%  '
$ws.Range("E30").Value = '  -0.75%  '
$ws.Range("D31").Value = "'0.0488"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.92%  '
$ws.Range("E32").Value = '  -0.62%  '
$ws.Range("D33").Value = "'3.22"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.15%  '
$ws.Range("E34").Value = '  -1.87%  '
$ws.Range("E35").Value = '  -0.47%  '
$ws.Range("D37").Value = '1.132.21'
$ws.Range("E37").Value = '  +1.70%  '
$ws.Range("E38").Value = '  -2.08%  '
$ws.Range("E39").Value = '  -2.23%  '
$ws.Range("D40").Value = "'0.0155"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.21%  '
$ws.Range("E41").Value = '  +0.12%  '
$ws.Range("E42").Value = '  -0.81%  '
$ws.Range("D43").Value = "'100.03"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.84%  '
$ws.Range("E44").Value = '  -1.32%  '
$ws.Range("E45").Value = '  -0.75%  '
$ws.Range("D46").Value = '1.764.85'
$ws.Range("E46").Value = '  -0.32%  '
$ws.Range("E47").Value = '  +1.99%  '
$ws.Range("D48").Value = "'55.18"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.42%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = "'0.0508"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.99%  '
$ws.Range("B50").Value = 'Mantle'
$ws.Range("C50").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D50").Value = "'0.418"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.05%  '
$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").Value = "'1.42"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +3.33%  '
